$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New commit row: "Juan camina y gira" / "Juan hace una animacion al caminar
# y gira en direccion al mouse" gets inserted as a new data row (row 8),
# pushed in right after the "Juan recupera vida" row (row 7). Duplicate the
# formatting of row 7 into row 8 first so the date/author/border/wrap styles
# match the rest of the table.
$ws.Range("A7:D7").Copy($ws.Range("A8:D8")) | Out-Null
$ws.Rows.Item(8).RowHeight = 45

$ws.Range("C8").Value = "Juan camina y gira "
$ws.Range("D8").Value = "Juan hace una animacion al caminar y gira en direccion al mouse "

# The existing "barra de vida" commit description (row 6) gains an extra
# sentence about a surprise.
$ws.Range("D6").Value = "Juan ya tiene una barra de vida arriba a la izquierda, y se creo un fuego para que este pierda 10% de su vida y muera y una sorpresa al suceder esto."

# Move the active selection to where it ended up after the edit.
$ws.Range("E6").Select() | Out-Null
